# Applies the "Add files via upload" revision to DocRequisitos.xlsx.
#
# Summary of the edit (from the OOXML diff):
#  - Sheet "Documento Geral de Recolha" (first sheet):
#      * C14 (R05.1 description) gains ",membros de seguranca" at the end.
#      * C27/C28 texts are unchanged (they just shift shared-string index,
#        which is an internal detail the object-model layer manages itself).
#      * C29 text is corrected: "resultados de provas" -> "resultados das provas".
#      * Row 30 (previously a blank placeholder row) is filled in with a new
#        requirement R21 ("Os atletas so competirao se passarem num teste de
#        dopping").
#      * Row 31 (previously a blank placeholder row) is filled in with a new
#        requirement R22 ("A competicao so tera inicio se todos os
#        funcionarios estiverem presentes").
#      * B32's number format is switched to match the dd/mm/yyyy h:mm format
#        used by the rows above it (cosmetic/style-only change).
#      * Column C is widened (more room for the longer descriptions).
#      * The active selection moves from G29 to C31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Update the existing "funcionario" requirement description (row 14)
#    to also mention "membros de seguranca".
# ---------------------------------------------------------------------
$ws.Range("C14").Value = "Cada funcionário tem um nome, um número de id e um cargo/função (juíz, camera, auxiliar,organizador,membros de segurança)"

# ---------------------------------------------------------------------
# 2. Fix the wording of the R20 requirement in row 29
#    ("de provas" -> "das provas").
# ---------------------------------------------------------------------
$ws.Range("C29").Value = "Só os juízes podem atualizar os resultados das provas durante a competição"

# ---------------------------------------------------------------------
# 3. Fill in row 30 with the new R21 requirement. The row already carries
#    the correct cell styles (same as the rows above), so plain value
#    assignment is enough.
# ---------------------------------------------------------------------
$ws.Range("A30").Value = "R21"
$ws.Range("B30").Value = 45566.8125
$ws.Range("C30").Value = "Os atletas só competirão se passarem num teste de dopping"
$ws.Range("F30").Value = "AC"
$ws.Range("G30").Value = "C"

# ---------------------------------------------------------------------
# 4. Fill in row 31 with the new R22 requirement. B31 currently uses a
#    different date style (m/d/yyyy) than the rows above (dd/mm/yyyy);
#    copy the number format from B29 first so it matches, then set values.
# ---------------------------------------------------------------------
$ws.Range("B29").Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A31").Value = "R22"
$ws.Range("B31").Value = 45566.814583333333
$ws.Range("C31").Value = "A competição só terá inicio se todos os funcionários estiverem presentes"
$ws.Range("F31").Value = "AC"
$ws.Range("G31").Value = "C"

# ---------------------------------------------------------------------
# 5. Row 32 keeps its (still blank) values, but its date cell B32 switches
#    to the same dd/mm/yyyy h:mm style used elsewhere in the table.
# ---------------------------------------------------------------------
$ws.Range("B29").Copy() | Out-Null
$ws.Range("B32").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------
# 6. Widen column C to fit the longer descriptions now in the table.
# ---------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 110

# ---------------------------------------------------------------------
# 7. Move the active selection to C31 (matching where the edit was made).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("C31").Select()

$excel.CutCopyMode = $false
